$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update client name in column A (rows 2-12) from "Adiko" to "Greba"
$ws.Range("A2:A12").Value = "Greba"

# Update the active selection to reflect the new selected range
$ws.Range("A2:A12").Select()
